# Se agrego el campo Nombres a la tabla de equivalencia del Modulo Gestionar Usuario
#
# This script replicates (on sheet "Clases de equialencia Nuevo U") a new
# equivalence-class block for the field "Nombres" in rows 11-13, mirroring
# the existing "Apellidos" block that lives in rows 8-10.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Copy the formatting (borders, fonts, alignment, fills) of the
#    "Apellidos" block (rows 8-10) down onto the new "Nombres" block
#    (rows 11-13) so the new rows look exactly the same way.
# ---------------------------------------------------------------------
$ws1.Range("B8:G10").Copy() | Out-Null
$ws1.Range("B11").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Helper color/font constants used by the little orange superscript-like
# suffix ("CEV<0x>"/"CENV<0x>" codes) that already exists elsewhere in
# the sheet (e.g. CEV<02>, CEV<03>, CENV<04>..CENV<06>).
# ---------------------------------------------------------------------
# FF993300 (ARGB) -> BGR integer used by the COM Font.Color property
$codeColor = 13209
$codeSize = 9
$codeFont = "Calibri"

function Set-CodeCell {
    param(
        $cell,
        [string]$code
    )
    $text = $code + " "
    $cell.Value2 = $text
    $len = $text.Length
    $chars = $cell.Characters($len, 1)
    $chars.Font.Size = $codeSize
    $chars.Font.Color = $codeColor
    $chars.Font.Name = $codeFont
}

# ---------------------------------------------------------------------
# 2) Row 11 - header row of the "Nombres" equivalence-class block
# ---------------------------------------------------------------------
$ws1.Range("B11").Value2 = "Nombres"
$ws1.Range("C11").Value2 = "Lógico"
$ws1.Range("D11").Value2 = "Nombres= caracteres alfanuméricos"
Set-CodeCell -cell $ws1.Range("E11") -code "CEV<04>"
$ws1.Range("F11").Value2 = "Nombres!= caracteres alfanuméricos"
Set-CodeCell -cell $ws1.Range("G11") -code "CENV<07>"

# ---------------------------------------------------------------------
# 3) Row 12
# ---------------------------------------------------------------------
$ws1.Range("C12").Value2 = "Valor"
$ws1.Range("D12").Value2 = "1 < Nombres<= 50"
Set-CodeCell -cell $ws1.Range("E12") -code "CEV<05>"
$ws1.Range("F12").Value2 = "Nombres<= 1"
Set-CodeCell -cell $ws1.Range("G12") -code "CENV<08>"

# ---------------------------------------------------------------------
# 4) Row 13
# ---------------------------------------------------------------------
$ws1.Range("F13").Value2 = "Nombres> 50"
Set-CodeCell -cell $ws1.Range("G13") -code "CENV<09>"

# ---------------------------------------------------------------------
# 5) Merge the cells the same way the "Apellidos" block (rows 8-10) is
#    merged, one row lower.
# ---------------------------------------------------------------------
$ws1.Range("B11:B13").Merge() | Out-Null
$ws1.Range("C12:C13").Merge() | Out-Null
$ws1.Range("D12:D13").Merge() | Out-Null
$ws1.Range("E12:E13").Merge() | Out-Null

# ---------------------------------------------------------------------
# 6) Move the active selection, as happened in the edited workbook.
# ---------------------------------------------------------------------
$ws1.Range("D21").Select() | Out-Null

Write-Output "Nombres equivalence-class block added"
